# Update the LR-pairs worksheet with recalculated TPM-based values.
# New data has only two sending clusters (FAPs, MuSCs) instead of three
# (ECs, FAPs, MuSCs) - so the workbook shrinks from 9 data rows to 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FAPs","Col2a1","Itgb1","ECs",1,0.3333333333333333,0.2021233333333333,0.60637,0.9764223557676824,0.9764223557676823,3,1,77.08952333333333,231.26857,0.2403816673726824,0.2403816673726824,15.58159142121111,140.2343227909,0.234714033939398,0.234714033939398),
    @("FAPs","Col2a1","Itgb1","FAPs",1,0.3333333333333333,0.2021233333333333,0.60637,0.9764223557676824,0.9764223557676823,3,1,101.5800373333333,304.740112,0.3167483425780597,0.3167483425780597,20.53169574593778,184.78526171344,0.309280162845578,0.3092801628455779),
    @("FAPs","Col2a1","Itgb1","MuSCs",1,0.3333333333333333,0.2021233333333333,0.60637,0.9764223557676824,0.9764223557676823,3,1,142.0267893333333,426.080368,0.4428699900492579,0.4428699900492579,28.70692808268445,258.36235274416,0.4324281589827065,0.4324281589827064),
    @("MuSCs","Col2a1","Itgb1","ECs",2,0.6666666666666666,0.004880666666666667,0.014642,0.02357764423231757,0.02357764423231757,3,1,77.08952333333333,231.26857,0.2403816673726824,0.2403816673726824,0.3762482668822222,3.38623440194,0.005667633433284407,0.005667633433284407),
    @("MuSCs","Col2a1","Itgb1","FAPs",2,0.6666666666666666,0.004880666666666667,0.014642,0.02357764423231757,0.02357764423231757,3,1,101.5800373333333,304.740112,0.3167483425780597,0.3167483425780597,0.4957783022115556,4.462004719904001,0.00746817973248174,0.007468179732481739),
    @("MuSCs","Col2a1","Itgb1","MuSCs",2,0.6666666666666666,0.004880666666666667,0.014642,0.02357764423231757,0.02357764423231757,3,1,142.0267893333333,426.080368,0.4428699900492579,0.4428699900492579,0.6931854164728889,6.238668748256001,0.01044183106655143,0.01044183106655143)
)

# Remove the three trailing rows (old "ECs" sending-cluster block, rows 8-10)
# so only 6 data rows remain (rows 2-7), matching the new TPM output.
$ws.Rows("8:10").Delete()

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value2 = $row[$c]
    }
}
